$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note text (cell A1) ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.85 = 10690.65 pesos`n✅ 10690.65 pesos = 2.84 = 966.59 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- tasas: update the rate cells N10, O10, N12, O12 ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 350.498
$ws2.Range("O10").Value = 3747.05
$ws2.Range("N12").Value = 3763.9
$ws2.Range("O12").Value = 340.311
